# Applies the addition of columns I (I0) and J (IF) to Sheet1,
# mirroring the header style of existing column H (IP) and filling in
# literal numeric values for rows 2-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting used by the other header cells (e.g. H1): bold
# font, thin border on all sides, centered horizontally, top-aligned
# vertically.
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous (thin)

# --- Data values for I2:J56 ---
$iValues = @(8,7,6,7,7,8,8,8,8,5,7,7,6,8,6,6,8,5,6,7,8,8,5,6,6,9,7,7,6,5,8,7,9,8,6,6,6,4,7,9,6,7,7,8,7,8,8,8,3,8,8,7,6,4,4)
$jValues = @(9,7,6,7,7,9,8,8,8,6,7,7,7,8,7,6,8,6,6,7,8,8,6,6,7,9,7,7,6,6,8,7,9,8,7,7,6,4,8,9,6,8,7,8,7,8,8,8,4,8,8,8,6,5,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]   # column I
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]  # column J
}
